$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded; insert it as row 4 (the
# existing rows 4-10 shift down to become rows 5-11, keeping their data).
$ws.Rows("4:4").Insert()

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 45092
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100108
$ws.Range("H4").Value = "Tropicales y subtropicales"
$ws.Range("I4").Value = 100108007
$ws.Range("J4").Value = "Coco"
$ws.Range("K4").Value = "Sin especificar"
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 150
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24333
$ws.Range("Q4").Value = "$/malla 20 unidades"
$ws.Range("R4").Value = "Perú"
$ws.Range("S4").Value = 1217
$ws.Range("T4").Value = 20
